$d = $word.ActiveDocument

# 1) "ithimbeni locwaningo" -> "ithimba locwaningo"
$d.Content.Find.Execute(
    "ithimbeni locwaningo ku-",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ithimba locwaningo ku-",
    2) | Out-Null

# 2) "uzizwa ukhululekile" -> "uzizwe ukhululekile"
$d.Content.Find.Execute(
    "uzizwa ukhululekile",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "uzizwe ukhululekile",
    2) | Out-Null

# 3) "kanti iMenenja yocwaningo nguZamakhanya Makhanya" -> "kanye neMenenja yocwaningo uZamakhanya Makhanya"
$d.Content.Find.Execute(
    "kanti iMenenja yocwaningo nguZamakhanya Makhanya",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "kanye neMenenja yocwaningo uZamakhanya Makhanya",
    2) | Out-Null

# 4) "noma okukukhathazayo mayelana namalungelo akho njengomhlanganyeli" ->
#    "noma kukhona okukukhathazayo mayelana namalungelo akho njengomhlanganyeli"
$d.Content.Find.Execute(
    "noma okukukhathazayo mayelana namalungelo akho njengomhlanganyeli",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "noma kukhona okukukhathazayo mayelana namalungelo akho njengomhlanganyeli",
    2) | Out-Null

# 5) "idokhumenti engenhla" -> "incwadi engenhla"
#    Use a tightly scoped Range + direct .Text assignment (rather than Find's
#    replace argument) so Word's smart-quote autoformat does not touch the
#    straight quotes around "EXIT" later in the same run/paragraph.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*idokhumenti engenhla*") {
        $r = $p.Range
        $r.Find.Execute("idokhumenti", $false) | Out-Null
        $r.Text = "incwadi"
        break
    }
}
